$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Save off the "before" values for rows 7, 8, 9 (only the columns that
# actually change: A, B, E, F, G, H, Q, R). The edit performs a cyclic
# rotation of these rows: old row 9 -> row 7, old row 7 -> row 8,
# old row 8 -> row 9.

$cols = @("A", "B", "E", "F", "G", "H", "Q", "R")

$row7 = @{}
$row8 = @{}
$row9 = @{}

foreach ($col in $cols) {
    $row7[$col] = $ws.Range("$col" + "7").Value2
    $row8[$col] = $ws.Range("$col" + "8").Value2
    $row9[$col] = $ws.Range("$col" + "9").Value2
}

foreach ($col in $cols) {
    $ws.Range("$col" + "7").Value2 = $row9[$col]
    $ws.Range("$col" + "8").Value2 = $row7[$col]
    $ws.Range("$col" + "9").Value2 = $row8[$col]
}
